$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken") and copy the header style from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:50:54.554100"
$ws.Range("F3").Value = "2021-10-05 10:50:54.554111"
$ws.Range("F4").Value = "2021-10-05 10:50:54.554114"
$ws.Range("F5").Value = "2021-10-05 10:50:54.554117"
$ws.Range("F6").Value = "2021-10-05 10:50:54.554120"
$ws.Range("F7").Value = "2021-10-05 10:50:54.554122"
$ws.Range("F8").Value = "2021-10-05 10:50:54.554125"
